$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 10:05"

$ws.Range("A6").Value = "Rusia"
$ws.Range("B6").Value = 232243
$ws.Range("C6").Value = 10899
$ws.Range("D6").Value = 43512
$ws.Range("E6").Value = 186615
$ws.Range("F6").Value = 2300
$ws.Range("G6").Value = 107
$ws.Range("H6").Value = 2116

$ws.Range("A7").Value = "Reino Unido"
$ws.Range("B7").Value = 223060
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 190651
$ws.Range("F7").Value = 1559
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 32065

$ws.Range("A28").Value = "Singapur"
$ws.Range("B28").Value = 24671
$ws.Range("C28").Value = 849
$ws.Range("D28").Value = 3225
$ws.Range("E28").Value = 21425
$ws.Range("F28").Value = 24
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 21

$ws.Range("A29").Value = "Bielorrusia"
$ws.Range("B29").Value = 23906
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 6531
$ws.Range("E29").Value = 17240
$ws.Range("F29").Value = 92
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 135

$ws.Range("D34").Value = 6131
$ws.Range("E34").Value = 9384

$ws.Range("B42").Value = 11350
$ws.Range("C42").Value = 264
$ws.Range("D42").Value = 2106
$ws.Range("E42").Value = 8493
$ws.Range("F42").Value = 31
$ws.Range("G42").Value = 25
$ws.Range("H42").Value = 751

$ws.Range("B66").Value = 3721
$ws.Range("C66").Value = 148
$ws.Range("D66").Value = 1250
$ws.Range("E66").Value = 2454

$ws.Range("B84").Value = 1746
$ws.Range("C84").Value = 5
$ws.Range("D84").Value = 777
$ws.Range("E84").Value = 908

$ws.Range("B89").Value = 1491
$ws.Range("C89").Value = 6
$ws.Range("D89").Value = 850
$ws.Range("E89").Value = 591

$ws.Range("A90").Value = "Eslovaquia"
$ws.Range("B90").Value = 1465
$ws.Range("C90").Value = 8
$ws.Range("D90").Value = 983
$ws.Range("E90").Value = 455
$ws.Range("F90").Value = 6
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 27

$ws.Range("A91").Value = "Eslovenia"
$ws.Range("B91").Value = 1460
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 256
$ws.Range("E91").Value = 1102
$ws.Range("F91").Value = 10
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 102

$ws.Range("A145").Value = "Nepal"
$ws.Range("B145").Value = 191
$ws.Range("C145").Value = 57
$ws.Range("D145").Value = 33
$ws.Range("E145").Value = 158
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0

$ws.Range("A146").Value = "Martinica"
$ws.Range("D146").Value = 91
$ws.Range("E146").Value = 82
$ws.Range("F146").Value = 4
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 14

$ws.Range("A147").Value = "Islas Feroe"
$ws.Range("B147").Value = 187
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 187
$ws.Range("E147").Value = 0
$ws.Range("F147").Value = 0

$ws.Range("A148").Value = "Madagascar"
$ws.Range("B148").Value = 186
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 101
$ws.Range("E148").Value = 85
$ws.Range("F148").Value = 1
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 0

$ws.Range("A149").Value = "Togo"
$ws.Range("B149").Value = 181
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 89
$ws.Range("E149").Value = 81
$ws.Range("H149").Value = 11

$ws.Range("A150").Value = "Birmania"
$ws.Range("B150").Value = 180
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 74
$ws.Range("E150").Value = 100
$ws.Range("H150").Value = 6

$ws.Range("A151").Value = "Suazilandia"
$ws.Range("B151").Value = 175
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 28
$ws.Range("E151").Value = 145
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 2

$ws.Range("A152").Value = "Sudan del Sur"
$ws.Range("B152").Value = 156
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 2
$ws.Range("E152").Value = 154
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 0

$ws.Range("A153").Value = "Guadalupe"
$ws.Range("B153").Value = 154
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 104
$ws.Range("E153").Value = 37
$ws.Range("F153").Value = 4
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 13

$ws.Range("A154").Value = "Gibraltar"
$ws.Range("B154").Value = 147
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 143
$ws.Range("E154").Value = 4
$ws.Range("H154").Value = 0

$ws.Range("A155").Value = "Guayana Francesa"
$ws.Range("B155").Value = 144
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 122
$ws.Range("E155").Value = 21
$ws.Range("H155").Value = 1

$ws.Range("A156").Value = "Republica de Africa Central"
$ws.Range("B156").Value = 143
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 10
$ws.Range("E156").Value = 133
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 0

$ws.Range("A157").Value = "Brunei"
$ws.Range("B157").Value = 141
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 134
$ws.Range("E157").Value = 6
$ws.Range("F157").Value = 2
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 1

$ws.Range("D177").Value = 42
$ws.Range("E177").Value = 3

$ws.Range("A215").Value = "San Bartolome"

$ws.Range("A216").Value = "Sahara Occidental"
